$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4154788255691528
$ws.Range("B1").Value = 0.8102923631668091
$ws.Range("C1").Value = 5.24896764755249
$ws.Range("D1").Value = 3.204691171646118
$ws.Range("E1").Value = 1.9486323595047
